$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poplar Samples")

# Shift existing column A data down by one row, matching the author's
# insertion of a new header row (A1 becomes the "kmerGWAS" header).
$ws.Rows("1:1").Insert()

# Header row
$ws.Range("A1").Value = 'kmerGWAS'
$ws.Range("D1").Value = 'GEMMA'
$ws.Range("E1").Value = 'coded sex'

# Column D (GEMMA fastq sample list) and E (coded sex) values
$ws.Range("D2").Value = 'SRR1820075_1.fastq SRR1820075_1.fastq 1'
$ws.Range("E2").Value = 'female'
$ws.Range("D3").Value = 'SRR1820078_1.fastq SRR1820078_1.fastq 1'
$ws.Range("E3").Value = 'female'
$ws.Range("D4").Value = 'SRR1820081_1.fastq SRR1820081_1.fastq 1'
$ws.Range("E4").Value = 'female'
$ws.Range("D5").Value = 'SRR1820090_1.fastq SRR1820090_1.fastq 1'
$ws.Range("E5").Value = 'female'
$ws.Range("D6").Value = 'SRR1820094_1.fastq SRR1820094_1.fastq 1'
$ws.Range("E6").Value = 'female'
$ws.Range("D7").Value = 'SRR1820099_1.fastq SRR1820099_1.fastq 1'
$ws.Range("E7").Value = 'female'
$ws.Range("D8").Value = 'SRR1820100_1.fastq SRR1820100_1.fastq 1'
$ws.Range("E8").Value = 'female'
$ws.Range("D9").Value = 'SRR1820101_1.fastq SRR1820101_1.fastq 1'
$ws.Range("E9").Value = 'female'
$ws.Range("D10").Value = 'SRR1820102_1.fastq SRR1820102_1.fastq 1'
$ws.Range("E10").Value = 'female'
$ws.Range("D11").Value = 'SRR1820103_1.fastq SRR1820103_1.fastq 1'
$ws.Range("E11").Value = 'female'
$ws.Range("D12").Value = 'SRR1820104_1.fastq SRR1820104_1.fastq 1'
$ws.Range("E12").Value = 'female'
$ws.Range("D13").Value = 'SRR1820110_1.fastq SRR1820110_1.fastq 1'
$ws.Range("E13").Value = 'female'
$ws.Range("D14").Value = 'SRR1820119_1.fastq SRR1820119_1.fastq 1'
$ws.Range("E14").Value = 'female'
$ws.Range("D15").Value = 'SRR1821206_1.fastq SRR1821206_1.fastq 1'
$ws.Range("E15").Value = 'female'
$ws.Range("D16").Value = 'SRR1821207_1.fastq SRR1821207_1.fastq 1'
$ws.Range("E16").Value = 'female'
$ws.Range("D17").Value = 'SRR1821208_1.fastq SRR1821208_1.fastq 1'
$ws.Range("E17").Value = 'female'
$ws.Range("D18").Value = 'SRR1821209_1.fastq SRR1821209_1.fastq 1'
$ws.Range("E18").Value = 'female'
$ws.Range("D19").Value = 'SRR1821210_1.fastq SRR1821210_1.fastq 1'
$ws.Range("E19").Value = 'female'
$ws.Range("D20").Value = 'SRR1821212_1.fastq SRR1821212_1.fastq 1'
$ws.Range("E20").Value = 'female'
$ws.Range("D21").Value = 'SRR1821213_1.fastq SRR1821213_1.fastq 1'
$ws.Range("E21").Value = 'female'
$ws.Range("D22").Value = 'SRR1821214_1.fastq SRR1821214_1.fastq 1'
$ws.Range("E22").Value = 'female'
$ws.Range("D23").Value = 'SRR1821215_1.fastq SRR1821215_1.fastq 1'
$ws.Range("E23").Value = 'female'
$ws.Range("D24").Value = 'SRR1821216_1.fastq SRR1821216_1.fastq 1'
$ws.Range("E24").Value = 'female'
$ws.Range("D25").Value = 'SRR1821223_1.fastq SRR1821223_1.fastq 1'
$ws.Range("E25").Value = 'female'
$ws.Range("D26").Value = 'SRR1821225_1.fastq SRR1821225_1.fastq 1'
$ws.Range("E26").Value = 'female'
$ws.Range("D27").Value = 'SRR1821232_1.fastq SRR1821232_1.fastq 1'
$ws.Range("E27").Value = 'female'
$ws.Range("D28").Value = 'SRR1821233_1.fastq SRR1821233_1.fastq 1'
$ws.Range("E28").Value = 'female'
$ws.Range("D29").Value = 'SRR1820107_1.fastq SRR1820107_1.fastq 1'
$ws.Range("E29").Value = 'female'
$ws.Range("D30").Value = 'SRR1820115_1.fastq SRR1820115_1.fastq 1'
$ws.Range("E30").Value = 'female'
$ws.Range("D31").Value = 'SRR1820117_1.fastq SRR1820117_1.fastq 1'
$ws.Range("E31").Value = 'female'
$ws.Range("D32").Value = 'SRR1821211_1.fastq SRR1821211_1.fastq 1'
$ws.Range("E32").Value = 'female'
$ws.Range("D33").Value = 'SRR1821227_1.fastq SRR1821227_1.fastq 1'
$ws.Range("E33").Value = 'female'
$ws.Range("D34").Value = 'SRR1821228_1.fastq SRR1821228_1.fastq 1'
$ws.Range("E34").Value = 'female'
$ws.Range("D35").Value = 'SRR1821230_1.fastq SRR1821230_1.fastq 1'
$ws.Range("E35").Value = 'female'
$ws.Range("D36").Value = 'SRR1821237_1.fastq SRR1821237_1.fastq 2'
$ws.Range("E36").Value = 'male'
$ws.Range("D37").Value = 'SRR1821238_1.fastq SRR1821238_1.fastq 2'
$ws.Range("E37").Value = 'male'
$ws.Range("D38").Value = 'SRR1821239_1.fastq SRR1821239_1.fastq 2'
$ws.Range("E38").Value = 'male'
$ws.Range("D39").Value = 'SRR1821240_1.fastq SRR1821240_1.fastq 2'
$ws.Range("E39").Value = 'male'
$ws.Range("D40").Value = 'SRR1821241_1.fastq SRR1821241_1.fastq 2'
$ws.Range("E40").Value = 'male'
$ws.Range("D41").Value = 'SRR1821244_1.fastq SRR1821244_1.fastq 2'
$ws.Range("E41").Value = 'male'
$ws.Range("D42").Value = 'SRR1821245_1.fastq SRR1821245_1.fastq 2'
$ws.Range("E42").Value = 'male'
$ws.Range("D43").Value = 'SRR1821250_1.fastq SRR1821250_1.fastq 2'
$ws.Range("E43").Value = 'male'
$ws.Range("D44").Value = 'SRR1821252_1.fastq SRR1821252_1.fastq 2'
$ws.Range("E44").Value = 'male'
$ws.Range("D45").Value = 'SRR1821253_1.fastq SRR1821253_1.fastq 2'
$ws.Range("E45").Value = 'male'
$ws.Range("D46").Value = 'SRR1821254_1.fastq SRR1821254_1.fastq 2'
$ws.Range("E46").Value = 'male'
$ws.Range("D47").Value = 'SRR1821235_1.fastq SRR1821235_1.fastq 2'
$ws.Range("E47").Value = 'male'
$ws.Range("D48").Value = 'SRR1821236_1.fastq SRR1821236_1.fastq 2'
$ws.Range("E48").Value = 'male'
$ws.Range("D49").Value = 'SRR1821246_1.fastq SRR1821246_1.fastq 2'
$ws.Range("E49").Value = 'male'
$ws.Range("D50").Value = 'SRR1821247_1.fastq SRR1821247_1.fastq 2'
$ws.Range("E50").Value = 'male'
$ws.Range("D51").Value = 'SRR1821248_1.fastq SRR1821248_1.fastq 2'
$ws.Range("E51").Value = 'male'
$ws.Range("D52").Value = 'SRR1821249_1.fastq SRR1821249_1.fastq 2'
$ws.Range("E52").Value = 'male'
$ws.Range("D53").Value = 'SRR1821251_1.fastq SRR1821251_1.fastq 2'
$ws.Range("E53").Value = 'male'

# Apply the same font styling used in column A (Menlo 11pt black) to
# the new column D cells, matching s="7" in the target worksheet.
$ws.Range("A2").Copy()
$ws.Range("D2:D53").PasteSpecial(-4122)

# Column D width, as set by the author for readability (best-fit to contents)
$ws.Columns("D").ColumnWidth = 46.95

# Restore the active selection to E2, matching the saved view state
$ws.Range("E2").Select()
